$d = $word.ActiveDocument

$replacements = @(
    @("29×46=", "26×25="),
    @("30×18=", "56×29="),
    @("89×52=", "65×81="),
    @("81×69=", "18×42="),
    @("33×82=", "61×20="),
    @("58×49=", "43×42="),
    @("33×45=", "58×62="),
    @("14×28=", "97×17="),
    @("25×11=", "90×94="),
    @("44×27=", "78×26="),
    @("29×40=", "94×17="),
    @("92×99=", "36×99="),
    @("51×17=", "62×53="),
    @("13×34=", "66×58="),
    @("74×28=", "43×26="),
    @("27×58=", "27×89="),
    @("49×85=", "47×67="),
    @("41×48=", "74×87="),
    @("58×48=", "50×34="),
    @("56×74=", "16×68="),
    @("82×45=", "24×62="),
    @("37×50=", "86×64="),
    @("99×36=", "74×70="),
    @("58×74=", "80×34="),
    @("31×75=", "25×78=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

$d.Save()
